# Generate Report for Handback
# - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
# - zh-cn / de-de sheets: populate "Latest Target File" (hyperlinked source .md),
#   "Latest Handback File" (xlf filename) and "Latest Handback DateTime" columns
#   for both data rows, and widen a few columns to fit the new content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# NOTE: the "Status" column cells on the zh-cn / de-de sheets (C2 / C3) point
# at the very same shared string ("Ready for handoff"); since that string's
# text is being swapped in place, those cells must be updated too so they
# keep sharing the (now re-worded) string - done further below alongside the
# rest of the per-language sheet edits.

# Widen the two status columns to fit the longer text (closest width the
# engine's character-grid rounding allows to the ~29.98 target).
$overview.Columns.Item(5).ColumnWidth = 29.14
$overview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column shares the "Ready for handoff" -> "Handed back: in sync with
# en-US" string with the Overview sheet.
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 - 355c1993...
$zhcn.Range("I2").Value = $zhcn.Range("A2").Value2
$linkI2 = $zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/355c1993-7e39-430e-b9bb-23a59c13e961.md")
$linkI2.TextToDisplay = "355c1993-7e39-430e-b9bb-23a59c13e961.md"
$zhcn.Range("J2").Value = "355c1993-7e39-430e-b9bb-23a59c13e961.d2e7200b482d807af86f85578b878dadf27a30ad.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-23 15:18:48"

# Row 3 - c730e75e...
$zhcn.Range("I3").Value = $zhcn.Range("A3").Value2
$linkI3 = $zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/c730e75e-4018-49cb-a52d-44a95de63869.md")
$linkI3.TextToDisplay = "c730e75e-4018-49cb-a52d-44a95de63869.md"
$zhcn.Range("J3").Value = "c730e75e-4018-49cb-a52d-44a95de63869.35360e632f992d3e636c17d6f43950c293aa2d72.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-23 15:18:48"

# Column widths: Status (C) and the two new text columns (I, J)
$zhcn.Columns.Item(3).ColumnWidth = 29.14
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column shares the "Ready for handoff" -> "Handed back: in sync with
# en-US" string with the Overview sheet.
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

# Row 2 - 355c1993...
$dede.Range("I2").Value = $dede.Range("A2").Value2
$linkDI2 = $dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/355c1993-7e39-430e-b9bb-23a59c13e961.md")
$linkDI2.TextToDisplay = "355c1993-7e39-430e-b9bb-23a59c13e961.md"
$dede.Range("J2").Value = "355c1993-7e39-430e-b9bb-23a59c13e961.d2e7200b482d807af86f85578b878dadf27a30ad.de-de.xlf"
$dede.Range("K2").Value = "2016-08-23 15:18:56"

# Row 3 - c730e75e...
$dede.Range("G3").Value = "c730e75e-4018-49cb-a52d-44a95de63869.35360e632f992d3e636c17d6f43950c293aa2d72.de-de.xlf"
$dede.Range("I3").Value = $dede.Range("A3").Value2
$linkDI3 = $dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a5caaae64a0155e5fa3e554a9e6c2b7c4310760/e2e/c730e75e-4018-49cb-a52d-44a95de63869.md")
$linkDI3.TextToDisplay = "c730e75e-4018-49cb-a52d-44a95de63869.md"
$dede.Range("J3").Value = "c730e75e-4018-49cb-a52d-44a95de63869.35360e632f992d3e636c17d6f43950c293aa2d72.de-de.xlf"
$dede.Range("K3").Value = "2016-08-23 15:18:56"

# Column widths: Status (C) and the two new text columns (I, J)
$dede.Columns.Item(3).ColumnWidth = 29.14
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
